$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "States testing with JDBC"
$ws.Range("B4").Value = "UNDEFINED"

$ws.Range("A5").Value = "States testing with JDBC"
$ws.Range("B5").Value = "PASSED"

$ws.Range("A6").Value = "States testing with JDBC"
$ws.Range("B6").Value = "PASSED"
